$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# The run "Your" / bookmark(_GoBack) / " Marks are " is merged back into a
# single run `"Your Marks are "`, dropping the stray _GoBack bookmark that
# used to sit between the two runs. A plain Find/Replace over the combined
# text re-merges the runs and removes the (zero-length) bookmark that used
# to live inside the replaced range.
$null = $d.Content.Find.Execute("Your Marks are", $true, $false, $false, $false, $false, $true, 1, $false, "Your Marks are", 2)

# --- Change 2 -------------------------------------------------------------
# Re-create the _GoBack bookmark as its own empty paragraph, inserted right
# after the "Call By Value : ... (default behaviour)" paragraph.
$found = $d.Content

$ok = $found.Find.Execute("Call By Value : When the parameter is passed as value to the called method (default behaviour)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$found.InsertParagraphAfter()

$newPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -eq ($found.End + 1)) {
        $newPara = $p
    }
}

$tr = $newPara.Range
$tr.InsertAfter("X")
$markRange = $d.Range($tr.Start, $tr.Start + 1)
$d.Bookmarks.Add("_GoBack", $markRange)
$markRange.Text = ""
